# After ClosedXML's "adjust to content" (AutoFit) algorithm was refactored,
# the computed column widths for this quote-prefixed-content sheet became
# narrower. Re-apply the new, narrower widths to columns A:B, C and D.
#
# Target widths, as persisted in the worksheet's <cols> collection (the
# "character width" unit stored in xl/worksheets/sheetN.xml):
#   A:B -> 2.996339
#   C   -> 9.282054
#   D   -> 11.710625
#
# Excel's ColumnWidth COM property uses the same character-width unit, but
# the host always reports/stores a column's persisted width as
# ColumnWidth + 5/6 (Excel pads every column by the equivalent of 5/6 of a
# character before writing the <col width="..."> attribute). Back that
# padding out of the desired, persisted widths so the value that actually
# lands in the file matches the target.
$padding = 5.0 / 6.0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 2.996339 - $padding
$ws.Columns.Item(2).ColumnWidth = 2.996339 - $padding
$ws.Columns.Item(3).ColumnWidth = 9.282054 - $padding
$ws.Columns.Item(4).ColumnWidth = 11.710625 - $padding
